$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)
$sh = $s.Shapes.Item("TextBox 10")
$tr = $sh.TextFrame.TextRange

# Replace the whole paragraph's text. This also removes the previous
# hyperlink ("aqui" -> https://.../form.dart) and its blue/underline
# styling, leaving plain black "Lato" 35pt text for the whole run -
# matching every remaining run's formatting in the target.
$tr.Text = "Acesse o código completo em form.dart"

# Re-touch each word/space range so the paragraph keeps the same
# per-word run boundaries as the authored slide (cosmetic split only -
# formatting is identical across all of them).
$ranges = @(
  @(1,6),   # "Acesse"
  @(7,3),   # " o "
  @(10,6),  # "código"
  @(16,1),  # " "
  @(17,8),  # "completo"
  @(25,1),  # " "
  @(26,2),  # "em"
  @(28,1),  # " "
  @(29,9)   # "form.dart"
)
foreach ($rg in $ranges) {
    $c = $tr.Characters($rg[0], $rg[1])
    $c.Font.Name = "Lato"
    $c.Font.Size = 35
    $c.Font.Color.RGB = 0
}

# The textbox uses <a:spAutoFit/>, so PowerPoint re-lays it out and
# resizes the shape to fit the new (now two-line) text.
$sh.Height = 94.11528
